# Weekly price-sheet update: a new week's prices (Valencia, Primera/Segunda)
# are inserted at the top of the data (rows 1186:1187), pushing all existing
# rows down by two. This also makes the two oldest rows, previously out of
# range, reappear at the bottom of the used range (rows 1254:1255) -
# Excel's row-insert shift handles that automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 1186, shifting rows 1186:1253
# down to 1188:1255 (and extending the sheet's dimension to A1:T1255).
$ws.Rows("1186:1187").Insert()

# New row 1186: Valencia / Primera
$ws.Range("A1186").Value = 3
$ws.Range("B1186").Value = "Femacal de La Calera"
$ws.Range("C1186").Value = "Coquimbo"
$ws.Range("D1186").Value = 45041
$ws.Range("E1186").Value = 5
$ws.Range("F1186").Value = "Fruta"
$ws.Range("G1186").Value = 100102
$ws.Range("H1186").Value = "Cítricos"
$ws.Range("I1186").Value = 100102005
$ws.Range("J1186").Value = "Naranja"
$ws.Range("K1186").Value = "Valencia"
$ws.Range("L1186").Value = "Primera"
$ws.Range("M1186").Value = 56
$ws.Range("N1186").Value = 14000
$ws.Range("O1186").Value = 14000
$ws.Range("P1186").Value = 14000
$ws.Range("Q1186").Value = "$/caja 15 kilos"
$ws.Range("R1186").Value = "Provincia de Quillota"
$ws.Range("S1186").Value = 933
$ws.Range("T1186").Value = 15

# New row 1187: Valencia / Segunda
$ws.Range("A1187").Value = 3
$ws.Range("B1187").Value = "Femacal de La Calera"
$ws.Range("C1187").Value = "Coquimbo"
$ws.Range("D1187").Value = 45041
$ws.Range("E1187").Value = 5
$ws.Range("F1187").Value = "Fruta"
$ws.Range("G1187").Value = 100102
$ws.Range("H1187").Value = "Cítricos"
$ws.Range("I1187").Value = 100102005
$ws.Range("J1187").Value = "Naranja"
$ws.Range("K1187").Value = "Valencia"
$ws.Range("L1187").Value = "Segunda"
$ws.Range("M1187").Value = 58
$ws.Range("N1187").Value = 12000
$ws.Range("O1187").Value = 12000
$ws.Range("P1187").Value = 12000
$ws.Range("Q1187").Value = "$/caja 15 kilos"
$ws.Range("R1187").Value = "Provincia de Quillota"
$ws.Range("S1187").Value = 800
$ws.Range("T1187").Value = 15
